$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.583.04"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").Value = "2.111.32"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.82%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.30%  "

$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5259"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4512"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.85%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09021"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "

$ws.Range("E11").Value = "  -0.60%  "

$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("D13").Value = "2.117.14"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.799"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.066"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "99.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001180"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.47%  "

$ws.Range("E18").Value = "  +0.90%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06707"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.33%  "

$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("E21").Value = "  +0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.318"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.49%  "

$ws.Range("D23").Value = "30.645.68"
$ws.Range("E23").Value = "  +0.37%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.384"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "

$ws.Range("D26").Value = "2.363.03"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.519"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.191"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.25%  "

$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.639"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.345"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.897"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.28%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02638"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06830"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2320"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6861"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.267"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.15%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.316"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.42%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6412"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.755"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.90%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000359"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.73%  "

$ws.Range("E49").Value = "  +0.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07285"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.59%  "
